# Fixed update to excel issue
# 1. Rename the "Requested quantity" headers to data-source-friendly names
$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2. Add the new "PO Forecast" sheet as the last tab
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (A2:D17)
$data = New-Object 'object[,]' 16,4
$data[0,0] = 45480.99999999999
$data[0,1] = 63
$data[0,2] = 19.96951041727689
$data[0,3] = 102.9768637564333
$data[1,0] = 45522.99999999999
$data[1,1] = 54
$data[1,2] = 12.60330858457754
$data[1,3] = 94.70032429275709
$data[2,0] = 45550.99999999999
$data[2,1] = 47
$data[2,2] = 4.85822061046464
$data[2,3] = 87.39174757597918
$data[3,0] = 45578.99999999999
$data[3,1] = 40
$data[3,2] = -1.921906575124653
$data[3,3] = 79.51168641887963
$data[4,0] = 45585.99999999999
$data[4,1] = 39
$data[4,2] = -2.006836898442605
$data[4,3] = 76.94240365085813
$data[5,0] = 45592.99999999999
$data[5,1] = 37
$data[5,2] = -3.895840411031066
$data[5,3] = 75.63244654030419
$data[6,0] = 45606.99999999999
$data[6,1] = 34
$data[6,2] = -3.358880466462753
$data[6,3] = 74.08746205926995
$data[7,0] = 45613.99999999999
$data[7,1] = 32
$data[7,2] = -9.37746211189895
$data[7,3] = 72.48348259559667
$data[8,0] = 45620.99999999999
$data[8,1] = 31
$data[8,2] = -10.59105352362335
$data[8,3] = 72.96751844099617
$data[9,0] = 45627.99999999999
$data[9,1] = 29
$data[9,2] = -16.73847823367074
$data[9,3] = 69.1593482919817
$data[10,0] = 45634.99999999999
$data[10,1] = 27
$data[10,2] = -15.555651361224
$data[10,3] = 67.61932016673201
$data[11,0] = 45641.99999999999
$data[11,1] = 26
$data[11,2] = -16.16252335206938
$data[11,3] = 66.13559561296762
$data[12,0] = 45648.99999999999
$data[12,1] = 24
$data[12,2] = -17.45334573939968
$data[12,3] = 63.44968177707997
$data[13,0] = 45655.99999999999
$data[13,1] = 22
$data[13,2] = -17.03284597789738
$data[13,3] = 63.07455500636515
$data[14,0] = 45662.99999999999
$data[14,1] = 21
$data[14,2] = -20.65613962552042
$data[14,3] = 63.40055096488815
$data[15,0] = 45669.99999999999
$data[15,1] = 19
$data[15,2] = -20.88675234751509
$data[15,3] = 65.19374144289088
$wsForecast.Range("A2:D17").Value = $data

# 3. Match the styling used on the existing sheets:
#    header row style (bold, centered, bordered) and date-formatted first column
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A17").PasteSpecial(-4122)

$excel.CutCopyMode = $false
